$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 0.1539904918020898
$ws.Cells.Item(2,2).Value = 0.02967996784325758
$ws.Cells.Item(2,3).Value = 5.188364509534735
$ws.Cells.Item(2,4).Value = [double]"2.121490321930391e-07"
$ws.Cells.Item(2,5).Value = 1.061

$ws.Cells.Item(3,1).Value = 0.4758799425243964
$ws.Cells.Item(3,2).Value = 0.05053636257183669
$ws.Cells.Item(3,3).Value = 9.416584777900074
$ws.Cells.Item(3,4).Value = [double]"4.659987381596779e-21"
$ws.Cells.Item(3,5).Value = 2.618

$ws.Cells.Item(4,1).Value = 0.1427894820662683
$ws.Cells.Item(4,2).Value = 0.02920515459780759
$ws.Cells.Item(4,3).Value = 4.889187680485258
$ws.Cells.Item(4,4).Value = [double]"1.012529269310844e-06"
$ws.Cells.Item(4,5).Value = 1.012

$ws.Cells.Item(5,1).Value = 0.1386594449149919
$ws.Cells.Item(5,2).Value = 0.04397037212300751
$ws.Cells.Item(5,3).Value = 3.153474447000149
$ws.Cells.Item(5,4).Value = 0.001613393527079192
$ws.Cells.Item(5,5).Value = 1.901

$ws.Cells.Item(6,1).Value = 0.09782385431104501
$ws.Cells.Item(6,2).Value = 0.02739470066044493
$ws.Cells.Item(6,3).Value = 3.570904297278648
$ws.Cells.Item(6,4).Value = 0.0003557508725838566
$ws.Cells.Item(6,5).Value = 1.017

$ws.Cells.Item(7,1).Value = 0.09852980910460823
$ws.Cells.Item(7,2).Value = 0.03394227281447775
$ws.Cells.Item(7,3).Value = 2.902864214283885
$ws.Cells.Item(7,4).Value = 0.003697669214550133
$ws.Cells.Item(7,5).Value = 1.176

$ws.Cells.Item(8,1).Value = 0.1285591507903313
$ws.Cells.Item(8,2).Value = 0.04506537660620899
$ws.Cells.Item(8,3).Value = 2.852725539469223
$ws.Cells.Item(8,4).Value = 0.004334604092332723
$ws.Cells.Item(8,5).Value = 2.22

$ws.Cells.Item(9,1).Value = 0.2480619298614547
$ws.Cells.Item(9,2).Value = 0.03814501256845883
$ws.Cells.Item(9,3).Value = 6.503128801340885
$ws.Cells.Item(9,4).Value = [double]"7.866638400616443e-11"
$ws.Cells.Item(9,5).Value = 1.455

$ws.Cells.Item(10,1).Value = 0.1406045920890998
$ws.Cells.Item(10,2).Value = 0.03898327176966881
$ws.Cells.Item(10,3).Value = 3.606793009059289
$ws.Cells.Item(10,4).Value = 0.0003100047591186715
$ws.Cells.Item(10,5).Value = 1.813

$ws.Cells.Item(11,1).Value = 0.1166872759966023
$ws.Cells.Item(11,2).Value = 0.04863551945342295
$ws.Cells.Item(11,3).Value = 2.399219280640169
$ws.Cells.Item(11,4).Value = 0.01643007231275663
$ws.Cells.Item(11,5).Value = 2.603

$ws.Cells.Item(12,1).Value = -0.1778421866078817
$ws.Cells.Item(12,2).Value = 0.05230021027414258
$ws.Cells.Item(12,3).Value = -3.400410546643777
$ws.Cells.Item(12,4).Value = 0.000672847470307431
$ws.Cells.Item(12,5).Value = 2.811

$ws.Cells.Item(13,1).Value = 0.1782436511497816
$ws.Cells.Item(13,2).Value = 0.03476419805854268
$ws.Cells.Item(13,3).Value = 5.127218837311319
$ws.Cells.Item(13,4).Value = [double]"2.940536496215887e-07"
$ws.Cells.Item(13,5).Value = 1.174

$ws.Cells.Item(14,1).Value = -1.410908757278547
$ws.Cells.Item(14,2).Value = 0.03352621097114366
$ws.Cells.Item(14,3).Value = -42.08375227647795
$ws.Cells.Item(14,4).Value = 0
